# "Added tests for template container and fixed typo"
#
# The "Template Container" controller block (rows 84-88 of Sheet1) had its
# "Integration Tests" column (C) left blank. Integration tests were added
# for it, so mark all five endpoint rows as covered ("Yes"), matching the
# same convention used by every other controller section on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

foreach ($r in 84..88) {
    $ws.Cells.Item($r, 3).Value = "Yes"
}

# Reflect where the user ended up / was looking after making the edit.
$ws.Range("C94").Select()
$excel.ActiveWindow.ScrollRow = 61
